$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.830.62'
$ws.Range("E2").Value = '  +3.23%  '

$ws.Range("D3").Value = '4.041.38'
$ws.Range("E3").Value = '  +3.06%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.30'
$ws.Range("E5").Value = '  -1.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.80'
$ws.Range("E6").Value = '  +2.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.622'
$ws.Range("E7").Value = '  +1.16%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.736'
$ws.Range("E9").Value = '  +0.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.174'
$ws.Range("E10").Value = '  +1.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000335'
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.79'
$ws.Range("E12").Value = '  +10.21%  '

$ws.Range("E13").Value = '  +5.04%  '

$ws.Range("D14").Value = '4.696.77'
$ws.Range("E14").Value = '  +3.21%  '

$ws.Range("D15").Value = '4.060.08'
$ws.Range("E15").Value = '  +2.59%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.22'
$ws.Range("E16").Value = '  +6.29%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.27'
$ws.Range("E17").Value = '  +2.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.21'
$ws.Range("E18").Value = '  -1.96%  '

$ws.Range("E19").Value = '  -1.67%  '

$ws.Range("D20").Value = '71.978.12'
$ws.Range("E20").Value = '  +3.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.64'
$ws.Range("E21").Value = '  +1.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '95.45'
$ws.Range("E22").Value = '  +8.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.51'
$ws.Range("E23").Value = '  +4.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.46'
$ws.Range("E24").Value = '  +0.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.19'
$ws.Range("E25").Value = '  +5.29%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.07'
$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  +4.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.00'
$ws.Range("E28").Value = '  +1.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.08'
$ws.Range("E29").Value = '  +8.57%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '710.48'
$ws.Range("E30").Value = '  +2.23%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.46'
$ws.Range("E31").Value = '  +1.70%  '

$ws.Range("E32").Value = '  +2.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.95'
$ws.Range("E33").Value = '  +16.97%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '68.08'
$ws.Range("E34").Value = '  -0.49%  '

$ws.Range("D35").Value = '0.0₃0906'
$ws.Range("E35").Value = '  +8.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.443'
$ws.Range("E36").Value = '  -0.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.68'
$ws.Range("E37").Value = '  +24.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.78'
$ws.Range("E38").Value = '  +0.92%  '

$ws.Range("E39").Value = '  +2.49%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.14%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.04%  '

$ws.Range("E42").Value = '  +0.60%  '

$ws.Range("E43").Value = '  +1.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  +0.80%  '

$ws.Range("E45").Value = '  +4.04%  '

$ws.Range("E46").Value = '  +2.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.18'
$ws.Range("E47").Value = '  +2.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000281'
$ws.Range("E48").Value = '  +24.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.07'
$ws.Range("E49").Value = '  +6.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.32'
$ws.Range("E50").Value = '  +0.42%  '

$ws.Range("D51").Value = '0.0₆0343'
$ws.Range("E51").Value = '  +0.76%  '
